# Bot5 GUI no había error del bot
# Update the "Cntdad feriados" parameter on the parametrosInicio sheet
# from 0 to 1, and move the selection down to B13 (matching the
# post-edit cursor position left behind when the workbook was saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parametrosInicio")
$ws.Activate()

$ws.Range("B12").Value = 1

$ws.Range("B13").Select()
